$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DND 3 V 0.3")

$ws.Cells.Item(3, 1).Value = "20/06/2024 05:44:56"
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = "IBREALEST"
$ws.Cells.Item(3, 4).Value = "Indiabulls Real Estate Limited"

$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "532832"
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(3, 6).Value = 7.14
$ws.Cells.Item(3, 7).Value = 146.7
$ws.Cells.Item(3, 8).Value = 29886307
